# Update res_bus vm_pu values for "case with 380 kV done"
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2
$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.026678170106423
$ws.Cells.Item(2, 4).Value = 1.037683182034894
$ws.Cells.Item(2, 5).Value = 1.047264022031575
$ws.Cells.Item(2, 6).Value = 1.051086584224076
$ws.Cells.Item(2, 9).Value = 1.033697897942123
$ws.Cells.Item(2, 10).Value = 1.031840273906417
$ws.Cells.Item(2, 11).Value = 1.040473367139956
$ws.Cells.Item(2, 12).Value = 1.050027146429824
$ws.Cells.Item(2, 13).Value = 1.053839061221582
$ws.Cells.Item(2, 14).Value = 1.014649001050638

# Row 3
$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.027546023956109
$ws.Cells.Item(3, 4).Value = 1.038378760310587
$ws.Cells.Item(3, 5).Value = 1.048179154123589
$ws.Cells.Item(3, 6).Value = 1.051982591987258
$ws.Cells.Item(3, 9).Value = 1.033851872131018
$ws.Cells.Item(3, 10).Value = 1.032348409877175
$ws.Cells.Item(3, 11).Value = 1.040979105975721
$ws.Cells.Item(3, 12).Value = 1.05075381958331
$ws.Cells.Item(3, 13).Value = 1.054547430280392
$ws.Cells.Item(3, 14).Value = 1.014817150790023

# Row 4
$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.028107824699541
$ws.Cells.Item(4, 4).Value = 1.038828653967148
$ws.Cells.Item(4, 5).Value = 1.048772324273933
$ws.Cells.Item(4, 6).Value = 1.052562951437551
$ws.Cells.Item(4, 9).Value = 1.033949646427925
$ws.Cells.Item(4, 10).Value = 1.032676816263043
$ws.Cells.Item(4, 11).Value = 1.041305486273215
$ws.Cells.Item(4, 12).Value = 1.051224412464149
$ws.Cells.Item(4, 13).Value = 1.055005739320718
$ws.Cells.Item(4, 14).Value = 1.014925807538744

# Row 5
$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.028344062202341
$ws.Cells.Item(5, 4).Value = 1.039017741621627
$ws.Cells.Item(5, 5).Value = 1.049021935249726
$ws.Cells.Item(5, 6).Value = 1.052807072344395
$ws.Cells.Item(5, 9).Value = 1.033990305255696
$ws.Cells.Item(5, 10).Value = 1.032814783458466
$ws.Cells.Item(5, 11).Value = 1.041442487776241
$ws.Cells.Item(5, 12).Value = 1.051422340939166
$ws.Cells.Item(5, 13).Value = 1.055198398273949
$ws.Cells.Item(5, 14).Value = 1.014971451159847

# Row 6
$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.028383730792455
$ws.Cells.Item(6, 4).Value = 1.039049487428473
$ws.Cells.Item(6, 5).Value = 1.049063860190121
$ws.Cells.Item(6, 6).Value = 1.052848069350773
$ws.Cells.Item(6, 9).Value = 1.033997105903552
$ws.Cells.Item(6, 10).Value = 1.032837943163376
$ws.Cells.Item(6, 11).Value = 1.041465478650295
$ws.Cells.Item(6, 12).Value = 1.051455579311017
$ws.Cells.Item(6, 13).Value = 1.055230745686824
$ws.Cells.Item(6, 14).Value = 1.014979112821131

# Row 7
$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.028110981096348
$ws.Cells.Item(7, 4).Value = 1.038831180754948
$ws.Cells.Item(7, 5).Value = 1.048775658638148
$ws.Cells.Item(7, 6).Value = 1.052566212853716
$ws.Cells.Item(7, 9).Value = 1.0339501914647
$ws.Cells.Item(7, 10).Value = 1.032678660161553
$ws.Cells.Item(7, 11).Value = 1.041307317716017
$ws.Cells.Item(7, 12).Value = 1.051227056837374
$ws.Cells.Item(7, 13).Value = 1.055008313696584
$ws.Cells.Item(7, 14).Value = 1.014926417571724

# Row 8
$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.026971414646228
$ws.Cells.Item(8, 4).Value = 1.03791829472082
$ws.Cells.Item(8, 5).Value = 1.047573083368782
$ws.Cells.Item(8, 6).Value = 1.051389272632291
$ws.Cells.Item(8, 9).Value = 1.033750318299588
$ws.Cells.Item(8, 10).Value = 1.032012081248352
$ws.Cells.Item(8, 11).Value = 1.040644462579458
$ws.Cells.Item(8, 12).Value = 1.050272648124637
$ws.Cells.Item(8, 13).Value = 1.054078467985657
$ws.Cells.Item(8, 14).Value = 1.014705858321318

# Row 9
$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.02496525758099
$ws.Cells.Item(9, 4).Value = 1.036308275948851
$ws.Cells.Item(9, 5).Value = 1.045461862394896
$ws.Cells.Item(9, 6).Value = 1.049319888057039
$ws.Cells.Item(9, 9).Value = 1.033383928597092
$ws.Cells.Item(9, 10).Value = 1.030834533788935
$ws.Cells.Item(9, 11).Value = 1.039469848418608
$ws.Cells.Item(9, 12).Value = 1.048593880181202
$ws.Cells.Item(9, 13).Value = 1.052439618222598
$ws.Cells.Item(9, 14).Value = 1.014316093323765

# Row 10
$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.023629175391414
$ws.Cells.Item(10, 4).Value = 1.035234089147519
$ws.Cells.Item(10, 5).Value = 1.044059766693781
$ws.Cells.Item(10, 6).Value = 1.047943444690054
$ws.Cells.Item(10, 9).Value = 1.033130179114347
$ws.Cells.Item(10, 10).Value = 1.030047583420762
$ws.Cells.Item(10, 11).Value = 1.038682426464901
$ws.Cells.Item(10, 12).Value = 1.047476813686427
$ws.Cells.Item(10, 13).Value = 1.051346906536454
$ws.Cells.Item(10, 14).Value = 1.014055525831848

# Row 11
$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.023050974127403
$ws.Cells.Item(11, 4).Value = 1.034768775397666
$ws.Cells.Item(11, 5).Value = 1.043453940477636
$ws.Cells.Item(11, 6).Value = 1.047348196089476
$ws.Cells.Item(11, 9).Value = 1.033018062138571
$ws.Cells.Item(11, 10).Value = 1.029706383559772
$ws.Cells.Item(11, 11).Value = 1.038340449316829
$ws.Cells.Item(11, 12).Value = 1.046993629222583
$ws.Cells.Item(11, 13).Value = 1.050873733598668
$ws.Cells.Item(11, 14).Value = 1.013942530264184

# Row 12
$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.022836255385384
$ws.Cells.Item(12, 4).Value = 1.034595911161725
$ws.Cells.Item(12, 5).Value = 1.043229105075953
$ws.Cells.Item(12, 6).Value = 1.047127210114748
$ws.Cells.Item(12, 9).Value = 1.03297608087865
$ws.Cells.Item(12, 10).Value = 1.029579580910678
$ws.Cells.Item(12, 11).Value = 1.038213272050876
$ws.Cells.Item(12, 12).Value = 1.046814231348283
$ws.Cells.Item(12, 13).Value = 1.050697974257687
$ws.Cells.Item(12, 14).Value = 1.013900533796892

# Row 13
$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.022882310959467
$ws.Cells.Item(13, 4).Value = 1.03463299227737
$ws.Cells.Item(13, 5).Value = 1.043277324164379
$ws.Cells.Item(13, 6).Value = 1.047174607087836
$ws.Cells.Item(13, 9).Value = 1.032985101194718
$ws.Cells.Item(13, 10).Value = 1.029606783471403
$ws.Cells.Item(13, 11).Value = 1.038240558852024
$ws.Cells.Item(13, 12).Value = 1.046852709234628
$ws.Cells.Item(13, 13).Value = 1.050735675294018
$ws.Cells.Item(13, 14).Value = 1.013909543301635

# Row 14
$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.023033224362856
$ws.Cells.Item(14, 4).Value = 1.034754486906871
$ws.Cells.Item(14, 5).Value = 1.043435351517537
$ws.Cells.Item(14, 6).Value = 1.047329926947091
$ws.Cells.Item(14, 9).Value = 1.033014598804813
$ws.Cells.Item(14, 10).Value = 1.029695903350063
$ws.Cells.Item(14, 11).Value = 1.038329939892301
$ws.Cells.Item(14, 12).Value = 1.04697879853239
$ws.Cells.Item(14, 13).Value = 1.050859205306266
$ws.Cells.Item(14, 14).Value = 1.013939059328517

# Row 15
$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.023126213878998
$ws.Cells.Item(15, 4).Value = 1.034829340362161
$ws.Cells.Item(15, 5).Value = 1.043532743328
$ws.Cells.Item(15, 6).Value = 1.047425640039061
$ws.Cells.Item(15, 9).Value = 1.033032728752817
$ws.Cells.Item(15, 10).Value = 1.029750804357058
$ws.Cells.Item(15, 11).Value = 1.038384990433302
$ws.Cells.Item(15, 12).Value = 1.047056496731763
$ws.Cells.Item(15, 13).Value = 1.050935316021185
$ws.Cells.Item(15, 14).Value = 1.013957241840914

# Row 16
$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.023667555784402
$ws.Cells.Item(16, 4).Value = 1.035264966755047
$ws.Cells.Item(16, 5).Value = 1.0441000007057
$ws.Cells.Item(16, 6).Value = 1.047982965538679
$ws.Cells.Item(16, 9).Value = 1.033137572783692
$ws.Cells.Item(16, 10).Value = 1.030070218442661
$ws.Cells.Item(16, 11).Value = 1.038705101015789
$ws.Cells.Item(16, 12).Value = 1.047508891962862
$ws.Cells.Item(16, 13).Value = 1.051378309130868
$ws.Cells.Item(16, 14).Value = 1.014063021465512

# Row 17
$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.024007214764158
$ws.Cells.Item(17, 4).Value = 1.035538175625496
$ws.Cells.Item(17, 5).Value = 1.04445617286385
$ws.Cells.Item(17, 6).Value = 1.048332765833838
$ws.Cells.Item(17, 9).Value = 1.033202739038348
$ws.Cells.Item(17, 10).Value = 1.030270460021531
$ws.Cells.Item(17, 11).Value = 1.038905626200298
$ws.Cells.Item(17, 12).Value = 1.047792805737626
$ws.Cells.Item(17, 13).Value = 1.051656182217829
$ws.Cells.Item(17, 14).Value = 1.014129329483945

# Row 18
$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.024205363970264
$ws.Cells.Item(18, 4).Value = 1.035697515897386
$ws.Cells.Item(18, 5).Value = 1.044664046466036
$ws.Cells.Item(18, 6).Value = 1.048536871693628
$ws.Cells.Item(18, 9).Value = 1.033240533061429
$ws.Cells.Item(18, 10).Value = 1.030387214525047
$ws.Cells.Item(18, 11).Value = 1.039022490818122
$ws.Cells.Item(18, 12).Value = 1.047958457218602
$ws.Cells.Item(18, 13).Value = 1.051818258713911
$ws.Cells.Item(18, 14).Value = 1.014167989582546

# Row 19
$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.024272933098003
$ws.Cells.Item(19, 4).Value = 1.035751843737049
$ws.Cells.Item(19, 5).Value = 1.044734947047055
$ws.Cells.Item(19, 6).Value = 1.048606478895404
$ws.Cells.Item(19, 9).Value = 1.033253383126029
$ws.Cells.Item(19, 10).Value = 1.030427017464288
$ws.Cells.Item(19, 11).Value = 1.039062321916111
$ws.Cells.Item(19, 12).Value = 1.048014948428176
$ws.Cells.Item(19, 13).Value = 1.051873522199304
$ws.Cells.Item(19, 14).Value = 1.014181168909906

# Row 20
$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.023970769264967
$ws.Cells.Item(20, 4).Value = 1.035508864732853
$ws.Cells.Item(20, 5).Value = 1.044417946082705
$ws.Cells.Item(20, 6).Value = 1.04829522798277
$ws.Cells.Item(20, 9).Value = 1.033195769688473
$ws.Cells.Item(20, 10).Value = 1.03024898043749
$ws.Cells.Item(20, 11).Value = 1.038884121907363
$ws.Cells.Item(20, 12).Value = 1.04776233935939
$ws.Cells.Item(20, 13).Value = 1.051626369280142
$ws.Cells.Item(20, 14).Value = 1.014122216938844

# Row 21
$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.02298878271081
$ws.Cells.Item(21, 4).Value = 1.034718710471173
$ws.Cells.Item(21, 5).Value = 1.043388810993515
$ws.Cells.Item(21, 6).Value = 1.047284185908615
$ws.Cells.Item(21, 9).Value = 1.033005921763503
$ws.Cells.Item(21, 10).Value = 1.029669661575825
$ws.Cells.Item(21, 11).Value = 1.038303623582064
$ws.Cells.Item(21, 12).Value = 1.046941666203462
$ws.Cells.Item(21, 13).Value = 1.050822828836269
$ws.Cells.Item(21, 14).Value = 1.013930368278126

# Row 22
$ws.Cells.Item(22, 2).Value = 1.019999999999999
$ws.Cells.Item(22, 3).Value = 1.022371664175121
$ws.Cells.Item(22, 4).Value = 1.034221759545174
$ws.Cells.Item(22, 5).Value = 1.042742884376694
$ws.Cells.Item(22, 6).Value = 1.046649174746125
$ws.Cells.Item(22, 9).Value = 1.03288461309395
$ws.Cells.Item(22, 10).Value = 1.029305040467
$ws.Cells.Item(22, 11).Value = 1.037937764017829
$ws.Cells.Item(22, 12).Value = 1.046426130375796
$ws.Cells.Item(22, 13).Value = 1.050317601087146
$ws.Cells.Item(22, 14).Value = 1.013809601634491

# Row 23
$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.022698781987018
$ws.Cells.Item(23, 4).Value = 1.034485216273454
$ws.Cells.Item(23, 5).Value = 1.043085194531261
$ws.Cells.Item(23, 6).Value = 1.046985742027939
$ws.Cells.Item(23, 9).Value = 1.032949105100726
$ws.Cells.Item(23, 10).Value = 1.029498368718148
$ws.Cells.Item(23, 11).Value = 1.03813179582916
$ws.Cells.Item(23, 12).Value = 1.046699382201565
$ws.Cells.Item(23, 13).Value = 1.050585432356936
$ws.Cells.Item(23, 14).Value = 1.013873635820266

# Row 24
$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.023987237318895
$ws.Cells.Item(24, 4).Value = 1.035522109118854
$ws.Cells.Item(24, 5).Value = 1.04443521873541
$ws.Cells.Item(24, 6).Value = 1.048312189495649
$ws.Cells.Item(24, 9).Value = 1.033198919506242
$ws.Cells.Item(24, 10).Value = 1.030258686269966
$ws.Cells.Item(24, 11).Value = 1.038893839075637
$ws.Cells.Item(24, 12).Value = 1.04777610565186
$ws.Cells.Item(24, 13).Value = 1.051639840470251
$ws.Cells.Item(24, 14).Value = 1.014125430842055

# Row 25
$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.025483663796957
$ws.Cells.Item(25, 4).Value = 1.036724659782712
$ws.Cells.Item(25, 5).Value = 1.046006721317005
$ws.Cells.Item(25, 6).Value = 1.049854326138321
$ws.Cells.Item(25, 9).Value = 1.033480325328398
$ws.Cells.Item(25, 10).Value = 1.031139300954141
$ws.Cells.Item(25, 11).Value = 1.039774285956119
$ws.Cells.Item(25, 12).Value = 1.049027515518167
$ws.Cells.Item(25, 13).Value = 1.052863331987181
$ws.Cells.Item(25, 14).Value = 1.014416986241469
